$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VATRIM")

$ws.Range("B2").Value = 137904.648850625
$ws.Range("B3").Value = 210458.463730442
$ws.Range("B4").Value = 210980.331756006
$ws.Range("B5").Value = 225210.561170175
$ws.Range("B6").Value = 146427.80203678
$ws.Range("B7").Value = 225062.905344404
$ws.Range("B8").Value = 224928.218442242
$ws.Range("B9").Value = 237117.07187975
$ws.Range("B10").Value = 157644.128386562
$ws.Range("B11").Value = 239530.243511537
$ws.Range("B12").Value = 239054.768312878
$ws.Range("B13").Value = 254294.598531707
$ws.Range("B14").Value = 168442.201770023
$ws.Range("B15").Value = 253771.419097488
$ws.Range("B16").Value = 253146.835079821
$ws.Range("B17").Value = 266163.02640468
$ws.Range("B18").Value = 178275.815981421
$ws.Range("B19").Value = 261247.990768206
$ws.Range("B20").Value = 257983.205635213
$ws.Range("B21").Value = 269648.511466152
$ws.Range("B22").Value = 172028.648307968
$ws.Range("B23").Value = 252607.990146722
$ws.Range("B24").Value = 250229.249114979
$ws.Range("B25").Value = 262995.235315836
$ws.Range("B26").Value = 171658.286377613
$ws.Range("B27").Value = 258612.323467944
$ws.Range("B28").Value = 258053.811110823
$ws.Range("B29").Value = 273274.892889381
$ws.Range("B30").Value = 176211.878336385
$ws.Range("B31").Value = 264579.675571494
$ws.Range("B32").Value = 264747.530374998
$ws.Range("B33").Value = 280287.854512593
$ws.Range("B34").Value = 187459.159705524
$ws.Range("B35").Value = 274916.287713791
$ws.Range("B36").Value = 272946.719365582
$ws.Range("B37").Value = 285409.227463171
$ws.Range("B38").Value = 187354.35846069
$ws.Range("B39").Value = 268079.566508662
$ws.Range("B40").Value = 262304.555021472
$ws.Range("B41").Value = 270591.251679394
$ws.Range("B42").Value = 176262.988592942
$ws.Range("B43").Value = 268424.615578619
$ws.Range("B44").Value = 266876.444841363
$ws.Range("B45").Value = 285382.931766757
$ws.Range("B46").Value = 179271.328487933
$ws.Range("B47").Value = 286422.385679443
$ws.Range("B48").Value = 287418.184746405
$ws.Range("B49").Value = 304882.863235226
$ws.Range("B50").Value = 191338.894130187
$ws.Range("B51").Value = 294318.447957792
$ws.Range("B52").Value = 289982.263771679
$ws.Range("B53").Value = 303392.712708047
$ws.Range("B54").Value = 176585.389154963
$ws.Range("B55").Value = 272232.85182302
$ws.Range("B56").Value = 268022.616864661
$ws.Range("B57").Value = 282300.072906288
$ws.Range("B58").Value = 161423.829949217
$ws.Range("B59").Value = 260457.531171536
$ws.Range("B60").Value = 260679.575216219
$ws.Range("B61").Value = 279747.748574731
$ws.Range("B62").Value = 165259.115938263
$ws.Range("B63").Value = 261643.315226561
$ws.Range("B64").Value = 255162.649327431
$ws.Range("B65").Value = 261138.996297799
$ws.Range("B66").Value = 155359.7141283
$ws.Range("B67").Value = 267749.772786758
$ws.Range("B68").Value = 268495.279737862
$ws.Range("B69").Value = 298778.326487255
$ws.Range("B70").Value = 185414.877275438
$ws.Range("B71").Value = 335908.496258
$ws.Range("B72").Value = 339475.03176023
$ws.Range("B73").Value = 360590.868087797
$ws.Range("B74").Value = 213525.162319122
$ws.Range("B75").Value = 345206.84445403
$ws.Range("B76").Value = 338649.549913408
$ws.Range("B77").Value = 358457.513279751
$ws.Range("B78").Value = 206371.271645858
$ws.Range("B79").Value = 343883.213439199
$ws.Range("B80").Value = 343528.820894054
$ws.Range("B81").Value = 368717.739774364
$ws.Range("B82").Value = 228022.528242696
$ws.Range("B83").Value = 375058.188932663
$ws.Range("B84").Value = 370765.914077571
$ws.Range("B85").Value = 389282.772247919
$ws.Range("B86").Value = 226300.111224402
$ws.Range("B87").Value = 371019.114262162
$ws.Range("B88").Value = 366277.50890152
$ws.Range("B89").Value = 391100.715725791
$ws.Range("B90").Value = 223888.619538789
$ws.Range("B91").Value = 382235.763584194
$ws.Range("B92").Value = 381144.834408552
$ws.Range("B93").Value = 406710.03732149
$ws.Range("B94").Value = 232115.205298561
$ws.Range("B95").Value = 390490.496211642
$ws.Range("B96").Value = 387354.843533586
$ws.Range("B97").Value = 412513.741715984
$ws.Range("B98").Value = 235310.28195268
$ws.Range("B99").Value = 398831.485340709
$ws.Range("B100").Value = 395726.072297189
$ws.Range("B101").Value = 422246.520517208
